# Add a new "2022-Q4" sheet (placed right after "总计") with its fund
# holdings data, and insert a corresponding new row into the "总计"
# summary sheet (pushing the existing 2021-Q4 / 2020-Q4 rows down).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q4" worksheet before the existing "2021-Q4"
#    sheet, so the final tab order is: 总计, 2022-Q4, 2021-Q4, 2020-Q4
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($beforeSheet)
$newSheet.Name = "2022-Q4"

# Header row (columns B..H), matching the style used on the other
# fund-holding sheets.
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"


# Data rows. Columns B, D, E, F, G hold text that happens to look
# numeric (fund code, fund size, position %, market value, ...) - the
# leading apostrophe forces Excel to store them as text, same as the
# source data (e.g. preserves "097.80" style trailing/leading zeros).
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'159855"
$newSheet.Range("C2").Value = "银华中证影视主题ETF"
$newSheet.Range("D2").Value = "'1.01"
$newSheet.Range("E2").Value = "'97.80"
$newSheet.Range("F2").Value = "'3.60"
$newSheet.Range("G2").Value = "'0.0364"
$newSheet.Range("H2").Value = 10

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'516620"
$newSheet.Range("C3").Value = "国泰中证影视主题ETF"
$newSheet.Range("D3").Value = "'0.71"
$newSheet.Range("E3").Value = "'98.01"
$newSheet.Range("F3").Value = "'3.82"
$newSheet.Range("G3").Value = "'0.0271"
$newSheet.Range("H3").Value = 10

# ---------------------------------------------------------------------
# 2) Insert a new row in "总计" for 2022-Q4, above the existing
#    2021-Q4 row (which, along with 2020-Q4, shifts down by one row).
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Range("A2:D2").Insert()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.06

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
